# Auto-generated: update FFXIV Leve profit calculation sheets with refreshed market prices
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 16111.111
$ws.Range("J21").Value = 16111.111
$ws.Range("L21").Value = 16111.111
$ws.Range("N21").Value = -17047.111
$ws.Range("H23").Value = 16111.111
$ws.Range("J23").Value = 16111.111
$ws.Range("L23").Value = 16111.111
$ws.Range("N23").Value = -16579.111
$ws.Range("H28").Value = 1222.2222
$ws.Range("I28").Value = 1285.8572
$ws.Range("J28").Value = 999.5
$ws.Range("K28").Value = 1285.8572
$ws.Range("L28").Value = 999.5
$ws.Range("M28").Value = -800.8571999999999
$ws.Range("N28").Value = -1969.5
$ws.Range("H29").Value = 575.125
$ws.Range("I29").Value = 575.125
$ws.Range("K29").Value = 1725.375
$ws.Range("M29").Value = -1444.375
$ws.Range("H58").Value = 11771.111
$ws.Range("I58").Value = 625
$ws.Range("J58").Value = 14955.714
$ws.Range("K58").Value = 1875
$ws.Range("L58").Value = 44867.142
$ws.Range("M58").Value = -1725
$ws.Range("N58").Value = -45167.142
$ws.Range("H62").Value = 2303.6667
$ws.Range("I62").Value = 1952.5
$ws.Range("J62").Value = 3006
$ws.Range("K62").Value = 1952.5
$ws.Range("L62").Value = 3006
$ws.Range("M62").Value = -1328.5
$ws.Range("N62").Value = -4254
$ws.Range("H65").Value = 2303.6667
$ws.Range("I65").Value = 1952.5
$ws.Range("J65").Value = 3006
$ws.Range("K65").Value = 9762.5
$ws.Range("L65").Value = 15030
$ws.Range("M65").Value = -6642.5
$ws.Range("N65").Value = -21270
$ws.Range("H80").Value = 373.85715
$ws.Range("I80").Value = 391.22223
$ws.Range("J80").Value = 342.6
$ws.Range("K80").Value = 1173.66669
$ws.Range("L80").Value = 1027.8
$ws.Range("M80").Value = -175.66669
$ws.Range("N80").Value = -3023.8
$ws.Range("H83").Value = 373.85715
$ws.Range("I83").Value = 391.22223
$ws.Range("J83").Value = 342.6
$ws.Range("K83").Value = 3521.00007
$ws.Range("L83").Value = 3083.4
$ws.Range("M83").Value = 1470.99993
$ws.Range("N83").Value = -13067.4
$ws.Range("H112").Value = 20409612
$ws.Range("J112").Value = 1497.9556
$ws.Range("L112").Value = 4493.8668
$ws.Range("N112").Value = -6709.8668
$ws.Range("H130").Value = 42780
$ws.Range("J130").Value = 42780
$ws.Range("L130").Value = 42780
$ws.Range("N130").Value = -52820
$ws.Range("H132").Value = 19883928
$ws.Range("I132").Value = 20695190
$ws.Range("J132").Value = 7999.5
$ws.Range("K132").Value = 62085570
$ws.Range("L132").Value = 23998.5
$ws.Range("M132").Value = -62083040
$ws.Range("N132").Value = -29058.5
$ws.Range("H137").Value = 4340.6304
$ws.Range("I137").Value = 3884.5757
$ws.Range("J137").Value = 5498.3076
$ws.Range("K137").Value = 11653.7271
$ws.Range("L137").Value = 16494.9228
$ws.Range("M137").Value = -9103.7271
$ws.Range("N137").Value = -21594.9228
$ws.Range("H138").Value = 4591.831
$ws.Range("I138").Value = 2739
$ws.Range("J138").Value = 4931.5166
$ws.Range("K138").Value = 8217
$ws.Range("L138").Value = 14794.5498
$ws.Range("M138").Value = -3077
$ws.Range("N138").Value = -25074.5498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4281.5884
$ws.Range("I74").Value = 4511.88
$ws.Range("J74").Value = 3641.889
$ws.Range("K74").Value = 4511.88
$ws.Range("L74").Value = 3641.889
$ws.Range("M74").Value = -3637.88
$ws.Range("N74").Value = -5389.889
$ws.Range("H77").Value = 4281.5884
$ws.Range("I77").Value = 4511.88
$ws.Range("J77").Value = 3641.889
$ws.Range("K77").Value = 22559.4
$ws.Range("L77").Value = 18209.445
$ws.Range("M77").Value = -18191.4
$ws.Range("N77").Value = -26945.445

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3383.5217
$ws.Range("I20").Value = 4453.6665
$ws.Range("J20").Value = 2695.5715
$ws.Range("K20").Value = 4453.6665
$ws.Range("L20").Value = 2695.5715
$ws.Range("M20").Value = -4206.6665
$ws.Range("N20").Value = -3189.5715
$ws.Range("H86").Value = 2168.2222
$ws.Range("I86").Value = 1916.6666
$ws.Range("J86").Value = 2671.3333
$ws.Range("K86").Value = 1916.6666
$ws.Range("L86").Value = 2671.3333
$ws.Range("M86").Value = -793.6666
$ws.Range("N86").Value = -4917.3333
$ws.Range("H89").Value = 2168.2222
$ws.Range("I89").Value = 1916.6666
$ws.Range("J89").Value = 2671.3333
$ws.Range("K89").Value = 9583.333000000001
$ws.Range("L89").Value = 13356.6665
$ws.Range("M89").Value = -3967.333000000001
$ws.Range("N89").Value = -24588.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6074.95
$ws.Range("I31").Value = 1850
$ws.Range("J31").Value = 7885.643
$ws.Range("K31").Value = 1850
$ws.Range("L31").Value = 7885.643
$ws.Range("M31").Value = -1555
$ws.Range("N31").Value = -8475.643
$ws.Range("H34").Value = 6074.95
$ws.Range("I34").Value = 1850
$ws.Range("J34").Value = 7885.643
$ws.Range("K34").Value = 1850
$ws.Range("L34").Value = 7885.643
$ws.Range("M34").Value = -1648
$ws.Range("N34").Value = -8289.643
$ws.Range("H141").Value = 32241.666
$ws.Range("J141").Value = 32241.666
$ws.Range("L141").Value = 32241.666
$ws.Range("N141").Value = -42601.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 928.7917
$ws.Range("I113").Value = 752.41174
$ws.Range("J113").Value = 1357.1428
$ws.Range("K113").Value = 2257.23522
$ws.Range("L113").Value = 4071.4284
$ws.Range("M113").Value = -87.23522000000003
$ws.Range("N113").Value = -8411.428400000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6256.393
$ws.Range("I70").Value = 5635.7334
$ws.Range("J70").Value = 6972.5386
$ws.Range("K70").Value = 5635.7334
$ws.Range("L70").Value = 6972.5386
$ws.Range("M70").Value = -5365.7334
$ws.Range("N70").Value = -7512.5386
$ws.Range("H73").Value = 6256.393
$ws.Range("I73").Value = 5635.7334
$ws.Range("J73").Value = 6972.5386
$ws.Range("K73").Value = 5635.7334
$ws.Range("L73").Value = 6972.5386
$ws.Range("M73").Value = -4699.7334
$ws.Range("N73").Value = -8844.5386
$ws.Range("H74").Value = 39650
$ws.Range("J74").Value = 39650
$ws.Range("L74").Value = 39650
$ws.Range("N74").Value = -41522
$ws.Range("H77").Value = 39650
$ws.Range("J77").Value = 39650
$ws.Range("L77").Value = 118950
$ws.Range("N77").Value = -128310
$ws.Range("H126").Value = 4070.1287
$ws.Range("I126").Value = 2859.282
$ws.Range("J126").Value = 5593.4517
$ws.Range("K126").Value = 8577.846000000001
$ws.Range("L126").Value = 16780.3551
$ws.Range("M126").Value = -6107.846000000001
$ws.Range("N126").Value = -21720.3551
$ws.Range("H132").Value = 3042.6128
$ws.Range("I132").Value = 1707.3636
$ws.Range("J132").Value = 3777
$ws.Range("K132").Value = 5122.0908
$ws.Range("L132").Value = 11331
$ws.Range("M132").Value = -2592.0908
$ws.Range("N132").Value = -16391

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 43205.09
$ws.Range("J80").Value = 43205.09
$ws.Range("L80").Value = 43205.09
$ws.Range("N80").Value = -45451.09
$ws.Range("H81").Value = 99999
$ws.Range("J81").Value = 99999
$ws.Range("L81").Value = 99999
$ws.Range("N81").Value = -101995
$ws.Range("H83").Value = 43205.09
$ws.Range("J83").Value = 43205.09
$ws.Range("L83").Value = 129615.27
$ws.Range("N83").Value = -140847.27
$ws.Range("H84").Value = 99999
$ws.Range("J84").Value = 99999
$ws.Range("L84").Value = 299997
$ws.Range("N84").Value = -309981
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
